$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------
# Sheet "Metadata": restructure columns, add new ones, rename headers,
# and move/rewrite the values for row 2.
# -------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Metadata")

# Extend the header formatting (bold / border / centered style that
# already lives on A1:G1) onto the new header cells H1:K1 so they match
# the existing header look.
$ws1.Range("A1").Copy()
$ws1.Range("H1:K1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 1 headers (final layout) ---
$ws1.Range("A1").Value = "Bank Name"
$ws1.Range("B1").Value = "Branch Name"
$ws1.Range("C1").Value = "Customer Name"
$ws1.Range("D1").Value = "Account Number"
$ws1.Range("E1").Value = "Customer ID"
$ws1.Range("F1").Value = "Account Type"
$ws1.Range("G1").Value = "IFSC Code"
$ws1.Range("H1").Value = "Opening Balance"
$ws1.Range("I1").Value = "Closing Balance"
$ws1.Range("J1").Value = "Period"
$ws1.Range("K1").Value = "Generated On"

# --- Row 2 data (final layout) ---
$ws1.Range("A2").Value = "Bank of Tomorrow Ltd."
$ws1.Range("B2").Value = "MG Road, Bangalore"
$ws1.Range("C2").Value = "Mr. Rajiv Sharma"
$ws1.Range("D2").Value = "'987654321012"
$ws1.Range("E2").Value = "'"
$ws1.Range("F2").Value = "'"
$ws1.Range("G2").Value = "BOTM0001234"
$ws1.Range("H2").Value = "₹1,25,000.00"
$ws1.Range("I2").Value = "'73,000.00"
$ws1.Range("J2").Value = "01-May-2024 to 31-May-2024"
$ws1.Range("K2").Value = "'"

# -------------------------------------------------------------------
# Sheet "Transactions": drop the "Balance" column (E) entirely and
# normalize a few date strings.
# -------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Transactions")

$ws2.Range("A2").Value = "'02-May-2024"
$ws2.Range("C2").Value = "₹15,000.00"

$ws2.Range("A3").Value = "'03-05-2024"

$ws2.Range("A5").Value = "'10-05-2024"

# Remove the whole "Balance" column (E) - header + all values - and
# shift the dimension back down to A1:D5.
$ws2.Range("E1:E5").Delete()

Write-Host "Edit complete"
